# Discharge_July16.xlsx edit: add "new depth" (cm) tables to both station
# sheets, converting each station's second (adjusted-velocity) table's D
# column from inches to centimeters (x2.54) and re-deriving segment/Q.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("stn3")
$ws2 = $wb.Worksheets.Item("stn4")

# ---------------------------------------------------------------------
# stn3 ("Station 3") — new "new depth" table at rows 35-49, derived from
# the existing second table at rows 19-31 (X stays the same, V is pasted
# as a value, D is converted cm = D(old)*2.54, segment/Q recomputed).
# ---------------------------------------------------------------------
$ws1.Range("A35").Value = "new depth"
$ws1.Range("A35").Font.Bold = $true

$ws1.Range("A36").Value = "X"
$ws1.Range("B36").Value = "V"
$ws1.Range("C36").Value = "D"
$ws1.Range("D36").Value = "segment"
$ws1.Range("E36").Value = "Q"
$ws1.Range("F36").Value = "Qtotal"

$ws1.Range("A37").Value = 0.55
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Formula = "=C19*2.54"
$ws1.Range("D37").Formula = "=A37"
$ws1.Range("F37").Formula = "=SUM(E37:E49)"

$ws1.Range("A38").Value = 0.6
$ws1.Range("B38").Value = 0.14872000000000002
$ws1.Range("C38").Formula = "=C20*2.54"
$ws1.Range("D38").Formula = "=(A38+(A39-A38)/2)"
$ws1.Range("E38").Formula = "=(D38-D37)*(B38)*C38"

$ws1.Range("A39").Value = 0.65
$ws1.Range("B39").Value = 0.18304
$ws1.Range("C39").Formula = "=C21*2.54"
$ws1.Range("D39").Formula = "=(A39+(A40-A39)/2)"
$ws1.Range("E39").Formula = "=(D39-D38)*(B39)*C39"

$ws1.Range("A40").Value = 0.7
$ws1.Range("B40").Value = 0.2288
$ws1.Range("C40").Formula = "=C22*2.54"
$ws1.Range("D40").Formula = "=(A40+(A41-A40)/2)"
$ws1.Range("E40").Formula = "=(D40-D39)*(B40)*C40"

$ws1.Range("A41").Value = 0.75
$ws1.Range("B41").Value = 0.28028000000000003
$ws1.Range("C41").Formula = "=C23*2.54"
$ws1.Range("D41").Formula = "=(A41+(A42-A41)/2)"
$ws1.Range("E41").Formula = "=(D41-D40)*(B41)*C41"

$ws1.Range("A42").Value = 0.8
$ws1.Range("B42").Value = 0.28600000000000003
$ws1.Range("C42").Formula = "=C24*2.54"
$ws1.Range("D42").Formula = "=(A42+(A43-A42)/2)"
$ws1.Range("E42").Formula = "=(D42-D41)*(B42)*C42"

$ws1.Range("A43").Value = 0.85
$ws1.Range("B43").Value = 0.30316
$ws1.Range("C43").Formula = "=C25*2.54"
$ws1.Range("D43").Formula = "=(A43+(A44-A43)/2)"
$ws1.Range("E43").Formula = "=(D43-D42)*(B43)*C43"

$ws1.Range("A44").Value = 0.9
$ws1.Range("B44").Value = 0.28600000000000003
$ws1.Range("C44").Formula = "=C26*2.54"
$ws1.Range("D44").Formula = "=(A44+(A45-A44)/2)"
$ws1.Range("E44").Formula = "=(D44-D43)*(B44)*C44"

$ws1.Range("A45").Value = 0.95
$ws1.Range("B45").Value = 0.20592000000000002
$ws1.Range("C45").Formula = "=C27*2.54"
$ws1.Range("D45").Formula = "=(A45+(A46-A45)/2)"
$ws1.Range("E45").Formula = "=(D45-D44)*(B45)*C45"

$ws1.Range("A46").Value = 1
$ws1.Range("B46").Value = 0.09724
$ws1.Range("C46").Formula = "=C28*2.54"
$ws1.Range("D46").Formula = "=(A46+(A47-A46)/2)"
$ws1.Range("E46").Formula = "=(D46-D45)*(B46)*C46"

$ws1.Range("A47").Value = 1.05
$ws1.Range("B47").Value = 0.0572
$ws1.Range("C47").Formula = "=C29*2.54"
$ws1.Range("D47").Formula = "=(A47+(A48-A47)/2)"
$ws1.Range("E47").Formula = "=(D47-D46)*(B47)*C47"

$ws1.Range("A48").Value = 1.1
$ws1.Range("B48").Value = 0.0286
$ws1.Range("C48").Formula = "=C30*2.54"
$ws1.Range("D48").Formula = "=(A48+(A49-A48)/2)"
$ws1.Range("E48").Formula = "=(D48-D47)*(B48)*C48"

$ws1.Range("A49").Value = 1.15
$ws1.Range("B49").Value = 0
$ws1.Range("C49").Formula = "=C31*2.54"
$ws1.Range("D49").Formula = "=(A49+(A50-A49)/2)"
$ws1.Range("E49").Formula = "=(D49-D48)*(B49)*C49"

# ---------------------------------------------------------------------
# stn4 ("Station 4") — new "new depth" table at rows 32-44, derived from
# the existing second table at rows 17-27.
# ---------------------------------------------------------------------
$ws2.Range("A32").Value = "new depth"
$ws2.Range("A32").Font.Bold = $true

$ws2.Range("A33").Value = "X"
$ws2.Range("B33").Value = "V"
$ws2.Range("C33").Value = "D"
$ws2.Range("D33").Value = "segment"
$ws2.Range("E33").Value = "Q"
$ws2.Range("F33").Value = "Qtotal"

$ws2.Range("A34").Value = 0.7
$ws2.Range("B34").Value = 0
$ws2.Range("C34").Formula = "=C17*2.54"
$ws2.Range("D34").Formula = "=A34"
$ws2.Range("F34").Formula = "=SUM(E34:E52)"

$ws2.Range("A35").Value = 0.75
$ws2.Range("B35").Value = 0.13155999999999998
$ws2.Range("C35").Formula = "=C18*2.54"
$ws2.Range("D35").Formula = "=(A35+(A36-A35)/2)"
$ws2.Range("E35").Formula = "=(D35-D34)*(B35)*C35"

$ws2.Range("A36").Value = 0.8
$ws2.Range("B36").Value = 0.3432
$ws2.Range("C36").Formula = "=C19*2.54"
$ws2.Range("D36").Formula = "=(A36+(A37-A36)/2)"
$ws2.Range("E36").Formula = "=(D36-D35)*(B36)*C36"

$ws2.Range("A37").Value = 0.85
$ws2.Range("B37").Value = 0.35464
$ws2.Range("C37").Formula = "=C20*2.54"
$ws2.Range("D37").Formula = "=(A37+(A38-A37)/2)"
$ws2.Range("E37").Formula = "=(D37-D36)*(B37)*C37"

$ws2.Range("A38").Value = 0.9
$ws2.Range("B38").Value = 0.37751999999999997
$ws2.Range("C38").Formula = "=C21*2.54"
$ws2.Range("D38").Formula = "=(A38+(A39-A38)/2)"
$ws2.Range("E38").Formula = "=(D38-D37)*(B38)*C38"

$ws2.Range("A39").Value = 0.95
$ws2.Range("B39").Value = 0.18304
$ws2.Range("C39").Formula = "=C22*2.54"
$ws2.Range("D39").Formula = "=(A39+(A40-A39)/2)"
$ws2.Range("E39").Formula = "=(D39-D38)*(B39)*C39"

$ws2.Range("A40").Value = 1
$ws2.Range("B40").Value = 0.19448
$ws2.Range("C40").Formula = "=C23*2.54"
$ws2.Range("D40").Formula = "=(A40+(A41-A40)/2)"
$ws2.Range("E40").Formula = "=(D40-D39)*(B40)*C40"

$ws2.Range("A41").Value = 1.05
$ws2.Range("B41").Value = 0.16016
$ws2.Range("C41").Formula = "=C24*2.54"
$ws2.Range("D41").Formula = "=(A41+(A42-A41)/2)"
$ws2.Range("E41").Formula = "=(D41-D40)*(B41)*C41"

$ws2.Range("A42").Value = 1.1
$ws2.Range("B42").Value = 0.14300000000000002
$ws2.Range("C42").Formula = "=C25*2.54"
$ws2.Range("D42").Formula = "=(A42+(A43-A42)/2)"
$ws2.Range("E42").Formula = "=(D42-D41)*(B42)*C42"

$ws2.Range("A43").Value = 1.15
$ws2.Range("B43").Value = 0.09724
$ws2.Range("C43").Formula = "=C26*2.54"
$ws2.Range("D43").Formula = "=(A43+(A44-A43)/2)"
$ws2.Range("E43").Formula = "=(D43-D42)*(B43)*C43"

$ws2.Range("A44").Value = 1.2
$ws2.Range("B44").Value = 0
$ws2.Range("C44").Formula = "=C27*2.54"
$ws2.Range("D44").Formula = "=(A44+(A45-A44)/2)"
$ws2.Range("E44").Formula = "=(D44-D43)*(B44)*C44"

# ---------------------------------------------------------------------
# View state: stn3 becomes the active/selected tab (was stn4), with its
# window scrolled down and a new selected cell; stn4 loses its former
# "tabSelected" flag and gets its own new selection.
# ---------------------------------------------------------------------
$ws2.Range("E31").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("G15").Select() | Out-Null
